$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-08 (row 21)
$ws.Range("B21").Value = 6240
$ws.Range("C21").Value = 989
$ws.Range("D21").Value = 5629268
$ws.Range("E21").Value = 902.126282051282
$ws.Range("F21").Value = 8.31452872765146
$ws.Range("G21").Value = 4.324894514767941
$ws.Range("H21").Value = 28.49358544176954
